$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.667.66"
$ws.Range("E2").Value = "  -4.24%  "

$ws.Range("D3").Value = "1.847.15"
$ws.Range("E3").Value = "  -3.63%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.44"
$ws.Range("E5").Value = "  -3.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4263"
$ws.Range("E7").Value = "  -6.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3640"
$ws.Range("E8").Value = "  -4.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.57"
$ws.Range("E9").Value = "  -4.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07231"
$ws.Range("E10").Value = "  -6.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8996"
$ws.Range("E11").Value = "  -7.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.63"
$ws.Range("E12").Value = "  -7.40%  "

$ws.Range("D13").Value = "1.828.47"
$ws.Range("E13").Value = "  -6.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.587"
$ws.Range("E14").Value = "  -5.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.363"
$ws.Range("E15").Value = "  -5.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06851"
$ws.Range("E16").Value = "  -2.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "77.63"
$ws.Range("E18").Value = "  -7.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008839"
$ws.Range("E19").Value = "  -6.64%  "

$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.40"
$ws.Range("E21").Value = "  -7.34%  "

$ws.Range("D22").Value = "27.632.67"
$ws.Range("E22").Value = "  -4.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.966"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.71"
$ws.Range("E24").Value = "  -3.09%  "

$ws.Range("D25").Value = "2.045.83"
$ws.Range("E25").Value = "  -4.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.051"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.59"
$ws.Range("E27").Value = "  -2.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.26"
$ws.Range("E28").Value = "  -4.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.340"
$ws.Range("E29").Value = "  -4.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "111.63"
$ws.Range("E30").Value = "  -5.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.763"
$ws.Range("E31").Value = "  -4.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08915"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7802"
$ws.Range("E33").Value = "  -10.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.506"
$ws.Range("E34").Value = "  -11.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.868"
$ws.Range("E35").Value = "  -4.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.084"
$ws.Range("E36").Value = "  -12.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05442"
$ws.Range("E38").Value = "  -4.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.095"
$ws.Range("E39").Value = "  -4.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.986"
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01927"
$ws.Range("E41").Value = "  -5.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5053"
$ws.Range("E42").Value = "  -7.92%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.773"
$ws.Range("E43").Value = "  -9.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1635"
$ws.Range("E44").Value = "  -6.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.276"
$ws.Range("E45").Value = "  -11.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06623"
$ws.Range("E46").Value = "  -4.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.74"
$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4722"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.33"
$ws.Range("E49").Value = "  -6.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("E51").Value = "  -6.56%  "
